# Auto-generated edit script: expand Artn-Gfra1 LR-pairs table with new TPM data
# (FAPs/MuSCs/Resolving-Mac sending clusters x ECs/FAPs/MuSCs target clusters)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tableData = @(
    @("FAPs", "Artn", "Gfra1", "ECs", 2.0, 0.6666666666666666, 0.736875, 2.210625, 0.5924430993326582, 0.5924430993326582, 1.0, 0.3333333333333333, 0.03998533333333333, 0.119956, 0.001814551768531471, 0.00181455176853147, 0.0294641925, 0.2651777325, 0.001075018673648341, 0.001075018673648341),
    @("FAPs", "Artn", "Gfra1", "FAPs", 2.0, 0.6666666666666666, 0.736875, 2.210625, 0.5924430993326582, 0.5924430993326582, 3.0, 1.0, 15.933008, 47.799024, 0.7230468132755195, 0.7230468132755195, 11.74063527, 105.66571743, 0.4283640950195505, 0.4283640950195505),
    @("FAPs", "Artn", "Gfra1", "MuSCs", 2.0, 0.6666666666666666, 0.736875, 2.210625, 0.5924430993326582, 0.5924430993326582, 3.0, 1.0, 6.062935333333333, 18.188806, 0.275138634955949, 0.275138634955949, 4.467625473749999, 40.20862926375, 0.1630039856394593, 0.1630039856394593),
    @("MuSCs", "Artn", "Gfra1", "ECs", 2.0, 0.6666666666666666, 0.3936963333333334, 1.181089, 0.3165295008188679, 0.3165295008188679, 1.0, 0.3333333333333333, 0.03998533333333333, 0.119956, 0.001814551768531471, 0.00181455176853147, 0.01574207912044445, 0.141678712084, 0.0005743591655032604, 0.0005743591655032604),
    @("MuSCs", "Artn", "Gfra1", "FAPs", 2.0, 0.6666666666666666, 0.3936963333333334, 1.181089, 0.3165295008188679, 0.3165295008188679, 3.0, 1.0, 15.933008, 47.799024, 0.7230468132755195, 0.7230468132755195, 6.272766828570668, 56.454901457136, 0.2288656468747734, 0.2288656468747734),
    @("MuSCs", "Artn", "Gfra1", "MuSCs", 2.0, 0.6666666666666666, 0.3936963333333334, 1.181089, 0.3165295008188679, 0.3165295008188679, 3.0, 1.0, 6.062935333333333, 18.188806, 0.275138634955949, 0.275138634955949, 2.386955409970445, 21.482598689734, 0.08708949477859128, 0.08708949477859128),
    @("Resolving-Mac", "Artn", "Gfra1", "ECs", 1.0, 0.3333333333333333, 0.113219, 0.339657, 0.09102739984847392, 0.09102739984847392, 1.0, 0.3333333333333333, 0.03998533333333333, 0.119956, 0.001814551768531471, 0.00181455176853147, 0.004527099454666666, 0.04074389509199999, 0.0001651739293798697, 0.0001651739293798697),
    @("Resolving-Mac", "Artn", "Gfra1", "FAPs", 1.0, 0.3333333333333333, 0.113219, 0.339657, 0.09102739984847392, 0.09102739984847392, 3.0, 1.0, 15.933008, 47.799024, 0.7230468132755195, 0.7230468132755195, 1.803919232752, 16.235273094768, 0.06581707138119558, 0.06581707138119558),
    @("Resolving-Mac", "Artn", "Gfra1", "MuSCs", 1.0, 0.3333333333333333, 0.113219, 0.339657, 0.09102739984847392, 0.09102739984847392, 3.0, 1.0, 6.062935333333333, 18.188806, 0.275138634955949, 0.275138634955949, 0.6864394755046667, 6.177955279541999, 0.02504515453789848, 0.02504515453789848)
)

for ($i = 0; $i -lt $tableData.Count; $i++) {
    $dataRow = $tableData[$i]
    $targetRow = $i + 2
    for ($j = 0; $j -lt $dataRow.Count; $j++) {
        $ws.Cells.Item($targetRow, $j + 1).Value = $dataRow[$j]
    }
}
